$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the RUN column text ("run") that was accidentally duplicated into
# rows 3 and 4 - only row 2 should carry it now.
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()

# Move the active selection to A2 (was J5)
$ws.Range("A2").Select()
